# Sinhala (si) translation pass for
# "Email 4-2 [TEMPLATE] Partner email - reminder to submit documents.docx"
#
# Strategy: every change in the source diff is a pure text swap inside an
# existing <w:t> run (no structural/paragraph changes). We therefore use
# Word's Find/Replace machinery to swap the English strings for their
# Sinhala counterparts. Most strings are unique (or all their occurrences
# get the same replacement) so a plain ReplaceAll is safe. Two strings
# ( " or " and "We're excited to see you at the upcoming " ) occur more
# times in the document than they should be replaced, so those are
# handled by walking occurrence-by-occurrence and only touching the
# specific ones that the diff changes.

$d = $word.ActiveDocument

$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-AllText($doc, $findText, $replaceText) {
    $doc.Content.Find.Execute(
        $findText, $true, $false, $false, $false, $false,
        $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll
    ) | Out-Null
}

# Replace the occurrence-th match of $findText (1-based, in document
# order) with $replaceText, leaving any other matches untouched.
function Replace-NthOccurrence($doc, $findText, $occurrence, $replaceText) {
    $rng = $doc.Content
    for ($i = 1; $i -le $occurrence; $i++) {
        $found = $rng.Find.Execute(
            $findText, $true, $false, $false, $false, $false,
            $true, $wdFindContinue, $false, "", 0
        )
        if (-not $found) {
            throw "Could not find occurrence $i of '$findText'"
        }
        if ($i -lt $occurrence) {
            $rng.Collapse(0)
        }
    }
    $rng.Text = $replaceText
}

# 1 & 2: "English" heading / link label (both occurrences -> same text)
Replace-AllText $d "English" "ඉංග්‍රීසී"

# 3: "Brief" (table header)
Replace-AllText $d "Brief" "කෙටි නෝට්ටුව"

# 4: brief description paragraph
Replace-AllText $d `
    "An email sent to partners in the target country who RSVPed yes but haven’t sent their documents to us. It will be sent via customer.io" `
    "දත්ත යැවීමේදී අපට ලියාපදිංචි කළ පාර්ශවිකයින්ට මෙන්ම අපට ලේඛන යවන්ට වුවමනාව දැක්කාත් නැති පාර්ශවිකයන්ට ඊ-මේල් යැවිය යුතුය. එම ලිපිය customer.io හරහා යවනු ලබයි"

# 5: "Target audience" header
Replace-AllText $d "Target audience" "ඉලක්ක ප්‍රේක්ෂකයා"

# 6: target audience description
Replace-AllText $d `
    "Invited partners who haven’t submitted their documents" `
    "ලේඛන නැවැත්වී ඇති ආරාධිත සහකාරයින්"

# 7 & 16: "Don't forget to send your documents" title (both emails)
Replace-AllText $d `
    "Don’t forget to send your documents" `
    "ඔබගේ ලේඛන යැවීමට අමතක නොකරන්න"

# 8: "Hi " greeting (first email only)
Replace-AllText $d "Hi " "ආයුබෝවන් "

# 9 & 20: "[insert list of documents required]" placeholder (both emails)
Replace-AllText $d `
    "[insert list of documents required]" `
    "[එම ලේඛනලේ ලැයිස්තුවක් ඇතුලත් කරන්න]"

# 10: "Please send a copy of these documents to your country manager, "
Replace-AllText $d `
    "Please send a copy of these documents to your country manager, " `
    "කරුණාකර මෙම ලේඛනවල පිටපතක් ඔබගේ රටේ කළමනාකරු, "

# 11 & 23: ", at " (both emails, same replacement)
Replace-AllText $d ", at " ", වෙත "

# 12 & 24: " or " -- occurs 3 times in the doc (2 should change, the
# "live chat ... or WhatsApp" one must stay). Replace 1st and 3rd only.
Replace-NthOccurrence $d " or " 1 " හරහා හෝ "
Replace-NthOccurrence $d " or " 2 " හරහා හෝ "

# 13: " (WhatsApp), so that we can make the necessary arrangements..."
Replace-AllText $d `
    " (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation." `
    " (WhatsApp) හරහා යවන්න, එවිට අපට නවාතැන් සහ ප්‍රවාහන පහසුකම් ඇතුළුව ඔබට අවශ්‍ය විධිවිධාන සලසන්න පුළුවන්."

# 14: "If you have any questions, please contact your country manager." (period)
Replace-AllText $d `
    "If you have any questions, please contact your country manager." `
    "ඔබට කිසියම් ප්‍රශ්නයක් ඇත්නම්, කරුණාකර ඔබේ රටේ කළමනාකරු අමතන්න."

# 15 & 26: "We look forward to seeing you there!" (both emails)
Replace-AllText $d `
    "We look forward to seeing you there!" `
    "අපි ඔබව එහි දැකීමට බලාපොරොත්තු වෙමු!"

# 17: "Dear " greeting (second email only)
Replace-AllText $d "Dear " "හිතවත් "

# 18: "We're excited to see you at the upcoming " -- occurs twice, only
# the second (the one in the second email, followed by ". '" smart quote)
# should change.
Replace-NthOccurrence $d "We’re excited to see you at the upcoming " 2 "ඉදිරියට පැවැත්වෙන "

# 19: ". '" (smart quote) -- unique in the doc
Replace-AllText $d ". ‘" " හිදී ඔබව දැක ගැනීමට අපට සතුටකි. ‘"

# 21: "Please reply to this email with a copy of these documents..."
Replace-AllText $d `
    "Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation." `
    "මෙම ඊ-මේල්ට ප්‍රතිචාර දී ඔබගේ ලේඛනවල පිටපතක් යවන්න, එවිට අපට ඔබට අවශ්‍ය විධිවිධාන සලස්වන්න පුළුවන්, නවාතැන් සහ ප්‍රවාහන පහසුකම් ඇතුළුව."

# 22: "If you have any questions, please contact your country manager, " (comma)
Replace-AllText $d `
    "If you have any questions, please contact your country manager, " `
    "ඔබට කිසියම් ප්‍රශ්නයක් ඇත්නම්, කරුණාකර ඔබගේ රටේ කළමනාකරු, "

# 25: " (WhatsApp)." (second email, closes the comment range)
Replace-AllText $d " (WhatsApp)." " (WhatsApp) හරහා අමතන්න."

# 27: comment text "choose either one"
# Find/Replace does not operate on comment ranges in this runtime, so
# assign the comment's Range.Text directly instead.
$d.Comments.Item(1).Range.Text = "ඉන් දෙකට එකක් තෝරන්න"
